$d = $word.ActiveDocument

# --- Step 1: insert the two new method-stage paragraphs, right after
# "2.3.3 running on test set" and right before "2.4 implementing in iOS" ---
$rng = $d.Content
$found = $rng.Find.Execute("2.3.3 running on test set", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter("`r2.3.4 Non maximum suppression on positive windows`r2.3.5 Hard negative finding (adding false positive windows to training set)")

# --- Step 2: rename "2.4 implementing in iOS" heading ---
$rng2 = $d.Content
$rng2.Find.Execute("2.4 implementing in iOS", $false, $false, $false, $false, $false, $true, 1, $false, "2.4 human versus model test", 2) | Out-Null

# --- Step 3: move the _GoBack bookmark from the end of the "2.1.2
# classifying ROI" paragraph to the end of the newly added
# "2.3.5 Hard negative finding ..." paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$rng3 = $d.Content
$rng3.Find.Execute("Hard negative finding (adding false positive windows to training set)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng3)
